$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '75.915.31'
$ws.Range('E2').Value = '  +1.28%  '

Set-TextValue 'D3' '2.931.51'
$ws.Range('E3').Value = '  +3.99%  '

$ws.Range('E4').Value = '  +0.07%  '

Set-TextValue 'D5' '202.51'
$ws.Range('E5').Value = '  +7.93%  '

Set-TextValue 'D6' '597.40'
$ws.Range('E6').Value = '  +0.61%  '

$ws.Range('E7').Value = '  -0.02%  '

Set-TextValue 'D8' '0.551'
$ws.Range('E8').Value = '  +0.32%  '

Set-TextValue 'D9' '0.197'
$ws.Range('E9').Value = '  +3.30%  '

Set-TextValue 'D10' '2.931.95'
$ws.Range('E10').Value = '  +4.07%  '

Set-TextValue 'D11' '0.437'
$ws.Range('E11').Value = '  +16.00%  '

$ws.Range('E12').Value = '  +0.75%  '

Set-TextValue 'D13' '4.94'
$ws.Range('E13').Value = '  +0.58%  '

Set-TextValue 'D14' '3.477.21'
$ws.Range('E14').Value = '  +4.34%  '

Set-TextValue 'D15' '28.08'
$ws.Range('E15').Value = '  +4.25%  '

Set-TextValue 'D16' '75.856.94'
$ws.Range('E16').Value = '  +1.29%  '

Set-TextValue 'D17' '0.0000190'
$ws.Range('E17').Value = '  +1.63%  '

Set-TextValue 'D18' '2.939.82'
$ws.Range('E18').Value = '  +4.42%  '

Set-TextValue 'D19' '13.21'
$ws.Range('E19').Value = '  +7.32%  '

Set-TextValue 'D20' '8.90'
$ws.Range('E20').Value = '  -1.82%  '

Set-TextValue 'D21' '373.55'
$ws.Range('E21').Value = '  -0.90%  '

Set-TextValue 'D22' '2.30'
$ws.Range('E22').Value = '  +1.63%  '

Set-TextValue 'D23' '4.30'
$ws.Range('E23').Value = '  +5.27%  '

Set-TextValue 'D24' '71.75'
$ws.Range('E24').Value = '  +1.29%  '

$ws.Range('E25').Value = '  +0.18%  '

Set-TextValue 'D26' '3.084.55'
$ws.Range('E26').Value = '  +4.52%  '

Set-TextValue 'D27' '4.30'
$ws.Range('E27').Value = '  +3.27%  '

Set-TextValue 'D28' '9.71'
$ws.Range('E28').Value = '  -0.45%  '

$ws.Range('E29').Value = '  +4.45%  '

$ws.Range('E30').Value = '  -0.69%  '

$ws.Range('E31').Value = '  +0.14%  '

Set-TextValue 'D32' '7.82'
$ws.Range('E32').Value = '  +1.87%  '

Set-TextValue 'D33' '502.02'
$ws.Range('E33').Value = '  -2.33%  '

$ws.Range('E34').Value = '  +2.87%  '

Set-TextValue 'D35' '0.999'
$ws.Range('E35').Value = '  +0.00%  '

$ws.Range('B36').Value = 'Cronos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D36' '0.113'
$ws.Range('E36').Value = '  +32.09%  '

Set-TextValue 'D37' '164.95'
$ws.Range('E37').Value = '  +0.78%  '

$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D38' '20.28'
$ws.Range('E38').Value = '  +1.62%  '

$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 'D39' '19.63'
$ws.Range('E39').Value = '  +1.37%  '

$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D40' '0.376'
$ws.Range('E40').Value = '  +9.70%  '

Set-TextValue 'D41' '0.113'
$ws.Range('E41').Value = '  -4.44%  '

Set-TextValue 'D42' '182.10'
$ws.Range('E42').Value = '  -2.26%  '

$ws.Range('E43').Value = '  -0.01%  '

Set-TextValue 'D44' '5.00'
$ws.Range('E44').Value = '  -0.18%  '

$ws.Range('E45').Value = '  -0.09%  '

Set-TextValue 'D46' '40.18'
$ws.Range('E46').Value = '  +0.21%  '

$ws.Range('E47').Value = '  -1.05%  '

Set-TextValue 'D48' '2.35'
$ws.Range('E48').Value = '  +0.76%  '

$ws.Range('E49').Value = '  +0.15%  '

Set-TextValue 'D50' '3.75'
$ws.Range('E50').Value = '  +0.65%  '

Set-TextValue 'D51' '22.54'
$ws.Range('E51').Value = '  +7.86%  '
